$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: Prototipos de interfaz de usuario (05/29/2025) ---
$ws.Range("A5").Value = 45806
$ws.Range("A5").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B5").Value = "Prototipos de interfaz de usuario"

# --- Header row: add "Autor" column header (C1), not centered like A1/B1 ---
$ws.Range("C1").Value = "Autor"

# --- "Autor" values (Ignacio Roldan interned here first) ---
$ws.Range("C2").Value = "Ignacio Roldan"
$ws.Range("C3").Value = "Ignacio Roldan"
$ws.Range("C4").Value = "Ignacio Roldan"
$ws.Range("C5").Value = "Ignacio Roldan"

# --- Row 7: date cell formatted but cleared (kept blank), bullet text in B ---
$ws.Range("A7").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("A7").ClearContents()
$ws.Range("B7").Value = "usuario cambiarlo por persona"

# --- Rows 8-12: correction bullet list (B column only) ---
$ws.Range("B8").Value = "visibilidad de los atributos, cambiarlos por privados"
$ws.Range("B9").Value = "lista de rutinas dentro de cliente"
$ws.Range("B10").Value = "lista de suscripción dentro de cliente"
$ws.Range("B11").Value = "lista de progresos en el cliente"
$ws.Range("B12").Value = "multiplicidad de ejercicio a ejercicio asignado(1 a 0..*)"

# --- Row 6: Correcciones: header (06/16/2025) ---
$ws.Range("A6").Value = 45824
$ws.Range("A6").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B6").Value = "Correcciones:"

# --- Row 13: creación de una clase DiaRutina (07/01/2025) ---
$ws.Range("A13").Value = 45839
$ws.Range("A13").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B13").Value = "creación de una clase DiaRutina"
$ws.Range("C13").Value = "Ignacio Roldan"

# --- Row 14: Diagramas de secuencia core (07/05/2025) ---
$ws.Range("A14").Value = 45843
$ws.Range("A14").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B14").Value = "Diagramas de secuencia core"
$ws.Range("C14").Value = "Ignacio Roldan"

# --- Row 15: Diagrama de clases (07/07/2025) ---
$ws.Range("A15").Value = 45845
$ws.Range("A15").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B15").Value = "Diagrama de clases"
$ws.Range("C15").Value = "Ignacio Roldan"

# --- Column widths (approximate best-fit sizing for new content) ---
$ws.Columns.Item(2).ColumnWidth = 26.333333333333332
$ws.Columns.Item(3).ColumnWidth = 11.833333333333334

# --- Selection moves to C4, matching the author's last-edited cell ---
$ws.Range("C4").Select() | Out-Null
